# Update the NATMI LR-pairs (App-Lrp1) sheet with recomputed TPM-based values.
# Ligand avg/total expr (G/H), their cluster-specificity (I/J), receptor avg/total
# expr (M/N), their cluster-specificity (O/P), and the derived edge weights
# (Q/R) and edge specificities (S/T) all change because the underlying per-cell
# TPM values were updated upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "G"=118.0346986666667; "H"=354.104096; "I"=0.2666057129183408; "J"=0.2666057129183408; "M"=3.795192333333334; "N"=11.385577; "O"=0.01044213755712683; "P"=0.01044213755712683; "Q"=447.9643834470436; "R"=4031.679451023392; "S"=0.002783933527809181; "T"=0.002783933527809181 }
    3 = @{ "G"=118.0346986666667; "H"=354.104096; "I"=0.2666057129183408; "J"=0.2666057129183408; "M"=243.3763986666667; "N"=730.1291960000001; "O"=0.6696287328350964; "P"=0.6696287328350964; "Q"=28726.85987919854; "R"=258541.7389127868; "S"=0.1785268457081061; "T"=0.1785268457081061 }
    4 = @{ "G"=118.0346986666667; "H"=354.104096; "I"=0.2666057129183408; "J"=0.2666057129183408; "M"=29.801371; "N"=89.404113; "O"=0.08199584844219236; "P"=0.08199584844219235; "Q"=3517.595845838538; "R"=31658.36261254684; "S"=0.02186056163027492; "T"=0.02186056163027492 }
    5 = @{ "G"=118.0346986666667; "H"=354.104096; "I"=0.2666057129183408; "J"=0.2666057129183408; "M"=86.47679266666667; "N"=259.430378; "O"=0.2379332811655844; "P"=0.2379332811655844; "Q"=10207.26216406981; "R"=91865.35947662829; "S"=0.06343437205215066; "T"=0.06343437205215066 }
    6 = @{ "I"=0.4881754016778185; "J"=0.4881754016778186; "M"=3.795192333333334; "N"=11.385577; "O"=0.01044213755712683; "P"=0.01044213755712683; "Q"=820.2569646120012; "R"=7382.31268150801; "S"=0.005097594696325427; "T"=0.005097594696325427 }
    7 = @{ "I"=0.4881754016778185; "J"=0.4881754016778186; "M"=243.3763986666667; "N"=730.1291960000001; "O"=0.6696287328350964; "P"=0.6696287328350964; "Q"=52601.07222370555; "R"=473409.65001335; "S"=0.3268962756267818; "T"=0.3268962756267819 }
    8 = @{ "I"=0.4881754016778185; "J"=0.4881754016778186; "M"=29.801371; "N"=89.404113; "O"=0.08199584844219236; "P"=0.08199584844219235; "Q"=6440.986377168969; "R"=57968.87739452071; "S"=0.04002835624918079; "T"=0.04002835624918079 }
    9 = @{ "I"=0.4881754016778185; "J"=0.4881754016778186; "M"=86.47679266666667; "N"=259.430378; "O"=0.2379332811655844; "P"=0.2379332811655844; "Q"=18690.27581003791; "R"=168212.4822903412; "S"=0.1161531751055305; "T"=0.1161531751055305 }
    10 = @{ "G"=45.876452; "H"=137.629356; "I"=0.1036214293744632; "J"=0.1036214293744632; "M"=3.795192333333334; "N"=11.385577; "O"=0.01044213755712683; "P"=0.01044213755712683; "Q"=174.1099589109347; "R"=1566.989630198412; "S"=0.001082029219394248; "T"=0.001082029219394248 }
    11 = @{ "G"=45.876452; "H"=137.629356; "I"=0.1036214293744632; "J"=0.1036214293744632; "M"=243.3763986666667; "N"=730.1291960000001; "O"=0.6696287328350964; "P"=0.6696287328350964; "Q"=11165.2456713642; "R"=100487.2110422778; "S"=0.06938788644658321; "T"=0.06938788644658322 }
    12 = @{ "G"=45.876452; "H"=137.629356; "I"=0.1036214293744632; "J"=0.1036214293744632; "M"=29.801371; "N"=89.404113; "O"=0.08199584844219236; "P"=0.08199584844219235; "Q"=1367.181166215692; "R"=12304.63049594123; "S"=0.008496527018351824; "T"=0.008496527018351822 }
    13 = @{ "G"=45.876452; "H"=137.629356; "I"=0.1036214293744632; "J"=0.1036214293744632; "M"=86.47679266666667; "N"=259.430378; "O"=0.2379332811655844; "P"=0.2379332811655844; "Q"=3967.248427886285; "R"=35705.23585097657; "S"=0.02465498669013389; "T"=0.02465498669013389 }
    14 = @{ "G"=62.68962833333333; "H"=188.068885; "I"=0.1415974560293775; "J"=0.1415974560293775; "M"=3.795192333333334; "N"=11.385577; "O"=0.01044213755712683; "P"=0.01044213755712683; "Q"=237.9191968301828; "R"=2141.272771471645; "S"=0.001478580113597978; "T"=0.001478580113597978 }
    15 = @{ "G"=62.68962833333333; "H"=188.068885; "I"=0.1415974560293775; "J"=0.1415974560293775; "M"=243.3763986666667; "N"=730.1291960000001; "O"=0.6696287328350964; "P"=0.6696287328350964; "Q"=15257.1759775185; "R"=137314.5837976665; "S"=0.0948177250536253; "T"=0.09481772505362532 }
    16 = @{ "G"=62.68962833333333; "H"=188.068885; "I"=0.1415974560293775; "J"=0.1415974560293775; "M"=29.801371; "N"=89.404113; "O"=0.08199584844219236; "P"=0.08199584844219235; "Q"=1868.236871813778; "R"=16814.131846324; "S"=0.01161040354438483; "T"=0.01161040354438483 }
    17 = @{ "G"=62.68962833333333; "H"=188.068885; "I"=0.1415974560293775; "J"=0.1415974560293775; "M"=86.47679266666667; "N"=259.430378; "O"=0.2379332811655844; "P"=0.2379332811655844; "Q"=5421.197991732059; "R"=48790.78192558853; "S"=0.03369074731776934; "T"=0.03369074731776935 }
}

foreach ($rowKey in $data.Keys) {
    $rowVals = $data[$rowKey]
    foreach ($col in $rowVals.Keys) {
        $cellAddr = "$col$rowKey"
        $ws.Range($cellAddr).Value = $rowVals[$col]
    }
}
